# changing velocity gradient calculation
# Rows 40-52 of the "IR profile" sheet used to step down by -0.1 (row 44 only)
# or -0.25 (rows 45-52) per 1.5mm; rows 2-43 were flat (=B$2 / =Cprev).
# The new profile instead starts the temperature drop two rows earlier
# (row 40, where the "temp drop point" marker now lives) and ramps down
# at a steady -0.2 per step all the way to row 52 (rows 53-60 keep their
# existing -1 per step formulas, their values simply follow the new base).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IR profile")

# Move the "temp drop point" marker from D44 to D40.
$ws.Range("D44").Value = ""
$ws.Range("D40").Value = "temp drop point"

# Rows 40-52: Outer (B) / Inner (C) now step down by 0.2 from the previous row.
For ($r = 40; $r -le 52; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 2).Formula = "=B" + $prev + "-0.2"
    $ws.Cells.Item($r, 3).Formula = "=C" + $prev + "-0.2"
}
